$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2/B3 password value: "India$321" -> "India$4321"
$ws.Range("B2").Value = "India`$4321"
$ws.Range("B3").Value = "India`$4321"

# Update the active selection to B3 (matches sheetView/selection in diff)
$ws.Range("B3").Select()
